$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 'AU-4 (1),AU-4'
$ws.Range("A3").Value = 'SC-5 (2),CM-6 b,SC-5'
$ws.Range("A4").Value = 'AU-12 (3),AU-7 b,AC-6 (8),AU-7 a,CM-5 (1),AC-6 (9),AU-8 b'
$ws.Range("A5").Value = 'AC-17 (9),CM-7 b,CM-6 b,AC-17 (1)'
$ws.Range("A8").Value = 'IA-2 (12),IA-2 (11)'
$ws.Range("A15").Value = 'IA-2,IA-8,AU-3 (1)'
$ws.Range("A17").Value = 'AU-12 a,MA-4 (1) (a),AU-3,AU-3 (1),AU-12 c'
$ws.Range("A19").Value = 'IA-5 (1) (b),CM-6 b,IA-5 (1) (a)'
$ws.Range("A21").Value = 'AC-12,MA-4 (7),SC-10,MA-4 e'
$ws.Range("A22").Value = 'CM-6 b,AU-12 a,AU-7 a,CM-5 (1),AU-3,MA-4 (1) (a),AU-3 (1),AU-14 (1),AU-7 (1),AU-6 (4)'
$ws.Range("A25").Value = 'AU-12 a,MA-4 (1) (a),AU-3,AU-3 (1),AU-12 c'
$ws.Range("A29").Value = 'SC-8 (1),SC-8 (2),SC-8'
$ws.Range("A31").Value = 'AU-12 a,AC-2 (4),MA-4 (1) (a),AU-3,AU-3 (1),AU-12 c'
$ws.Range("A34").Value = 'AC-11 a,AC-11 b'
$ws.Range("A45").Value = 'AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 b,AC-8 a'
$ws.Range("A55").Value = 'AC-17 (2),SC-8'
$ws.Range("A65").Value = 'IA-2 (2),CM-6 b'
$ws.Range("A67").Value = 'AU-12 a,MA-4 (1) (a),AU-3,AU-3 (1),AU-12 c'
$ws.Range("A69").Value = 'AU-12 (3),CM-6 b,AU-7 b,AU-12 a,AU-7 a,CM-5 (1),AU-12 c,AU-8 b'
$ws.Range("A77").Value = 'AU-12 a,AC-2 (4),MA-4 (1) (a),AU-3,AU-3 (1),AU-12 c'
$ws.Range("A80").Value = 'IA-2 (2),IA-2 (1),IA-2 (3),IA-2 (4)'
$ws.Range("A86").Value = 'AU-12 a,MA-4 (1) (a),AU-3,AU-3 (1),AU-12 c'
$ws.Range("A88").Value = 'AC-6 (9),AC-2 (4),AU-12 c,CM-5 (1)'
$ws.Range("A89").Value = 'IA-2,IA-2 (5),IA-2 (2),IA-2 (4),IA-2 (3)'
$ws.Range("A90").Value = 'IA-2 (12),IA-2 (11)'
$ws.Range("A97").Value = 'AU-8 (1) (a),AU-8 b,AU-8 (1) (b)'
$ws.Range("A102").Value = 'AU-12 a,MA-4 (1) (a),AU-3,AU-3 (1),AU-12 c'
$ws.Range("A111").Value = 'AU-5 a,AU-5 b'
$ws.Range("A119").Value = 'AU-12 a,MA-4 (1) (a),AU-3,AU-3 (1),AU-12 c'
$ws.Range("A124").Value = 'AU-12 a,MA-4 (1) (a),AU-3,AU-3 (1),AU-12 c'
$ws.Range("A125").Value = 'CM-7 a,AC-18 (1)'
$ws.Range("A128").Value = 'IA-5 (1) (c),CM-6 b,CM-7 a'
$ws.Range("A136").Value = 'AC-11 b,AC-11 (1)'
$ws.Range("A139").Value = 'CM-3 (5),SI-6 b,SI-6 d'
$ws.Range("A148").Value = 'AU-12 a,MA-4 (1) (a),AU-3,AU-3 (1),AU-12 c,AU-14 (1)'
$ws.Range("A157").Value = 'AU-12 a,MA-4 (1) (a),AU-3,AU-3 (1),AU-12 c'
$ws.Range("A159").Value = 'AC-17 (2),SC-8'
